$d = $word.ActiveDocument

$d.Content.Find.Execute("917÷7=131, 0", $true, $false, $false, $false, $false, $true, 1, $false, "866÷4=216, 2", 2) | Out-Null
$d.Content.Find.Execute("468÷4=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "816÷7=116, 4", 2) | Out-Null
$d.Content.Find.Execute("183÷4=45, 3", $true, $false, $false, $false, $false, $true, 1, $false, "687÷3=229, 0", 2) | Out-Null
$d.Content.Find.Execute("329÷8=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "936÷7=133, 5", 2) | Out-Null
$d.Content.Find.Execute("665÷2=332, 1", $true, $false, $false, $false, $false, $true, 1, $false, "754÷2=377, 0", 2) | Out-Null
$d.Content.Find.Execute("237÷2=118, 1", $true, $false, $false, $false, $false, $true, 1, $false, "285÷7=40, 5", 2) | Out-Null
$d.Content.Find.Execute("823÷9=91, 4", $true, $false, $false, $false, $false, $true, 1, $false, "159÷4=39, 3", 2) | Out-Null
$d.Content.Find.Execute("231÷3=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "375÷6=62, 3", 2) | Out-Null
$d.Content.Find.Execute("574÷9=63, 7", $true, $false, $false, $false, $false, $true, 1, $false, "177÷9=19, 6", 2) | Out-Null
$d.Content.Find.Execute("204÷8=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "440÷7=62, 6", 2) | Out-Null
$d.Content.Find.Execute("245÷7=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "934÷6=155, 4", 2) | Out-Null
$d.Content.Find.Execute("992÷5=198, 2", $true, $false, $false, $false, $false, $true, 1, $false, "486÷2=243, 0", 2) | Out-Null
$d.Content.Find.Execute("736÷6=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "588÷8=73, 4", 2) | Out-Null
$d.Content.Find.Execute("229÷8=28, 5", $true, $false, $false, $false, $false, $true, 1, $false, "870÷7=124, 2", 2) | Out-Null
$d.Content.Find.Execute("485÷6=80, 5", $true, $false, $false, $false, $false, $true, 1, $false, "108÷2=54, 0", 2) | Out-Null
$d.Content.Find.Execute("139÷2=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "680÷4=170, 0", 2) | Out-Null
$d.Content.Find.Execute("755÷5=151, 0", $true, $false, $false, $false, $false, $true, 1, $false, "717÷4=179, 1", 2) | Out-Null
$d.Content.Find.Execute("257÷9=28, 5", $true, $false, $false, $false, $false, $true, 1, $false, "160÷4=40, 0", 2) | Out-Null
$d.Content.Find.Execute("986÷7=140, 6", $true, $false, $false, $false, $false, $true, 1, $false, "904÷6=150, 4", 2) | Out-Null
$d.Content.Find.Execute("411÷5=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "826÷9=91, 7", 2) | Out-Null
$d.Content.Find.Execute("716÷4=179, 0", $true, $false, $false, $false, $false, $true, 1, $false, "811÷2=405, 1", 2) | Out-Null
$d.Content.Find.Execute("136÷7=19, 3", $true, $false, $false, $false, $false, $true, 1, $false, "259÷4=64, 3", 2) | Out-Null
$d.Content.Find.Execute("202÷7=28, 6", $true, $false, $false, $false, $false, $true, 1, $false, "991÷3=330, 1", 2) | Out-Null
$d.Content.Find.Execute("921÷2=460, 1", $true, $false, $false, $false, $false, $true, 1, $false, "374÷4=93, 2", 2) | Out-Null
$d.Content.Find.Execute("143÷4=35, 3", $true, $false, $false, $false, $false, $true, 1, $false, "255÷4=63, 3", 2) | Out-Null
